$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 44; every existing
# record from row 44 down to row 179 shifts down by one row, and the
# record that used to be on the last row (179) becomes the new last
# row (180).

# Make sure the newly-created last row has the same date number format
# as the rest of column D before we populate it.
$ws.Range("D180").NumberFormat = $ws.Range("D179").NumberFormat

# Shift rows 44..179 down to 45..180, working from the bottom up so we
# never overwrite a source row before it has been copied.
for ($r = 179; $r -ge 44; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Populate the newly inserted row 44 with the new observation.
$ws.Range("A44").Value2 = 7
$ws.Range("B44").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value2 = "Ñuble"
$ws.Range("D44").Value2 = 44623
$ws.Range("E44").Value2 = 16
$ws.Range("F44").Value2 = 100112017
$ws.Range("G44").Value2 = "Apio"
$ws.Range("H44").Value2 = "Americana (o)"
$ws.Range("I44").Value2 = "Primera"
$ws.Range("J44").Value2 = 60
$ws.Range("K44").Value2 = 8000
$ws.Range("L44").Value2 = 9000
$ws.Range("M44").Value2 = 8500
$ws.Range("N44").Value2 = "$/docena de matas"
$ws.Range("O44").Value2 = "Provincia del Elquí"
$ws.Range("P44").Value2 = 1417
$ws.Range("Q44").Value2 = 6
$ws.Range("R44").Value2 = "Hortaliza"
